$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(300724004056, 'https://www.naukri.com/job-listings-data-analyst-data-science-sr-trainer-talentify-hyderabad-pune-bengaluru-6-to-11-years-300724004056', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Analyst- Data Science  Sr Trainer_300724004056.xlsx')
    ,@(290724500690, 'https://www.naukri.com/job-listings-analytics-data-science-and-iot-lead-hexaware-technologies-ltd-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-5-to-9-years-290724500690', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Analytics Data science and IOT Lead_290724500690.xlsx')
    ,@(290724500689, 'https://www.naukri.com/job-listings-analytics-data-science-and-iot-lead-hexaware-technologies-ltd-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-5-to-9-years-290724500689', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Analytics Data science and IOT Lead_290724500689.xlsx')
    ,@(140324012002, 'https://www.naukri.com/job-listings-data-science-manager-axtria-pune-gurugram-bengaluru-8-to-12-years-140324012002', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Manager_140324012002.xlsx')
    ,@(170724900002, 'https://www.naukri.com/job-listings-s-c-global-network-ai-cg-s-consultant-data-science-accenture-solutions-pvt-ltd-pune-5-to-10-years-170724900002', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/S&C Global Network - AI - CG&S - Consultant Data Science_170724900002.xlsx')
    ,@(240524012199, 'https://www.naukri.com/job-listings-hiring-data-science-engineer-lead-architect-xoriant-hyderabad-pune-mumbai-all-areas-5-to-10-years-240524012199', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Hiring Data Science Engineer - Lead - Architect !_240524012199.xlsx')
    ,@(240724909755, 'https://www.naukri.com/job-listings-data-science-manager-4bell-technology-mumbai-indore-pune-8-to-13-years-240724909755', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Manager_240724909755.xlsx')
    ,@(160724500785, 'https://www.naukri.com/job-listings-manager-data-science-genai-tmt-vertical-trendence-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-4-to-8-years-160724500785', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Manager - Data Science (GenAI) - TMT Vertical_160724500785.xlsx')
    ,@(20724502537, 'https://www.naukri.com/job-listings-senior-manager-data-science-and-analytics-metafolks-pune-5-to-10-years-020724502537', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Senior Manager - Data Science and Analytics_20724502537.xlsx')
    ,@(70323501288, 'https://www.naukri.com/job-listings-senior-analyst-data-science-navit-software-solutions-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-5-to-8-years-070323501288', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Senior Analyst Data Science_70323501288.xlsx')
    ,@(120724004528, 'https://www.naukri.com/job-listings-data-science-engineer-gen-ai-is-required-perm-role-ibm-pune-chennai-bengaluru-5-to-10-years-120724004528', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Engineer- GEN AI is required-PERM ROLE_120724004528.xlsx')
    ,@(20724010574, 'https://www.naukri.com/job-listings-data-science-manager-axa-global-business-services-pune-bengaluru-8-to-13-years-020724010574', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Manager_20724010574.xlsx')
    ,@(50724905401, 'https://www.naukri.com/job-listings-data-science-practitioner-accenture-solutions-pvt-ltd-pune-5-to-10-years-050724905401', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Practitioner_50724905401.xlsx')
    ,@(261023501559, 'https://www.naukri.com/job-listings-data-science-engineer-cognisure-pune-bengaluru-7-to-10-years-261023501559', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Engineer_261023501559.xlsx')
    ,@(231122500002, 'https://www.naukri.com/job-listings-bigdata-data-science-engineer-quantta-analytics-pvt-ltd-pune-4-to-8-years-231122500002', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/BigData-Data Science Engineer_231122500002.xlsx')
    ,@(40624501335, 'https://www.naukri.com/job-listings-data-science-nextcomm-corporation-india-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-1-to-8-years-040624501335', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science_40624501335.xlsx')
    ,@(30724010023, 'https://www.naukri.com/job-listings-faculty-data-science-analytics-prowess-publishing-hyderabad-pune-3-to-8-years-030724010023', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Faculty - Data Science & Analytics_30724010023.xlsx')
    ,@(290724502209, 'https://www.naukri.com/job-listings-data-scientist-agilite-global-solutions-company-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-5-to-10-years-290724502209', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Scientist_290724502209.xlsx')
    ,@(290724501408, 'https://www.naukri.com/job-listings-data-scientist-norstella-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-5-to-10-years-290724501408', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Scientist_290724501408.xlsx')
    ,@(230724013592, 'https://www.naukri.com/job-listings-senior-data-scientist-nvidia-intel-ai-ltimindtree-pune-chennai-bengaluru-5-to-10-years-230724013592', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Senior Data Scientist - Nvidia Intel AI_230724013592.xlsx')
    ,@(60924012222, 'https://www.naukri.com/job-listings-data-science-consultant-optimization-supply-chain-and-manufacturing-zs-pune-bengaluru-6-to-9-years-060924012222', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Consultant (Optimization - Supply Chain and Manufacturing_60924012222.xlsx')
    ,@(280824916175, 'https://www.naukri.com/job-listings-data-science-engineer-cybage-software-pvt-ltd-pune-8-to-12-years-280824916175', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Engineer_280824916175.xlsx')
    ,@(260824913880, 'https://www.naukri.com/job-listings-data-science-professional-capgemini-technology-services-india-limited-hyderabad-pune-bengaluru-6-to-9-years-260824913880', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/DATA SCIENCE Professional_260824913880.xlsx')
    ,@(210824500836, 'https://www.naukri.com/job-listings-analytics-data-science-and-iot-lead-hexaware-technologies-ltd-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-6-to-9-years-210824500836', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Analytics Data science and IOT Lead_210824500836.xlsx')
    ,@(10824501631, 'https://www.naukri.com/job-listings-analytics-data-science-and-iot-lead-hexaware-technologies-ltd-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-5-to-9-years-010824501631', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Analytics Data science and IOT Lead_10824501631.xlsx')
    ,@(80824500328, 'https://www.naukri.com/job-listings-senior-platform-automation-data-science-engineer-mastercard-pune-4-to-8-years-080824500328', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Senior Platform Automation & Data Science Engineer_80824500328.xlsx')
    ,@(220824503831, 'https://www.naukri.com/job-listings-lead-data-science-engineer-wolters-kluwer-medknow-publications-pune-5-to-10-years-220824503831', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Lead Data Science Engineer_220824503831.xlsx')
    ,@(230824502289, 'https://www.naukri.com/job-listings-manager-data-science-tredence-kolkata-pune-chennai-gurugram-bengaluru-8-to-13-years-230824502289', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Manager - Data Science_230824502289.xlsx')
    ,@(250824000014, 'https://www.naukri.com/job-listings-data-science-engineer-r-studio-integrated-personnel-services-pune-chennai-bengaluru-4-to-8-years-250824000014', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Engineer (R Studio)_250824000014.xlsx')
    ,@(140824913740, 'https://www.naukri.com/job-listings-lead-data-science-nitor-infotech-pvt-ltd-pune-6-to-9-years-140824913740', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Lead Data Science_140824913740.xlsx')
    ,@(100824005069, 'https://www.naukri.com/job-listings-lead-data-science-nitor-pune-6-to-9-years-100824005069', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Lead Data Science_100824005069.xlsx')
    ,@(281022502443, 'https://www.naukri.com/job-listings-data-science-professional-extentia-information-technology-pune-4-to-9-years-281022502443', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Science Professional_281022502443.xlsx')
    ,@(271122001064, 'https://www.naukri.com/job-listings-data-scientist-consultant-aarons-visions-hyderabad-pune-bengaluru-6-to-11-years-271122001064', '/home/ravindra/git_repos/monitoring_ds_jobs/data/scraped_data/Data Scientist Consultant_271122001064.xlsx')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = [double]$data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}